# Weekly update: two new Granada price records observed at the market are
# inserted at the top of the data block (rows 3-4), pushing the previously
# existing rows 3-10 down to rows 5-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 3 (existing rows 3-10
# shift down to rows 5-12, formatting/formulas moving along with them).
$ws.Rows("3:4").Insert()

# ---- New row 3 ----
$ws.Cells.Item(3, 1).Value = 6
$ws.Cells.Item(3, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(3, 3).Value = "Metropolitana"
$ws.Cells.Item(3, 4).Value = 44662
$ws.Cells.Item(3, 5).Value = 13
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100104
$ws.Cells.Item(3, 8).Value = "Frutos de pepita"
$ws.Cells.Item(3, 9).Value = 100104001
$ws.Cells.Item(3, 10).Value = "Granada"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 45
$ws.Cells.Item(3, 14).Value = 18000
$ws.Cells.Item(3, 15).Value = 18000
$ws.Cells.Item(3, 16).Value = 18000
$ws.Cells.Item(3, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(3, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 19).Value = 1000
$ws.Cells.Item(3, 20).Value = 18

# ---- New row 4 ----
$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(4, 3).Value = "Metropolitana"
$ws.Cells.Item(4, 4).Value = 44662
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100104
$ws.Cells.Item(4, 8).Value = "Frutos de pepita"
$ws.Cells.Item(4, 9).Value = 100104001
$ws.Cells.Item(4, 10).Value = "Granada"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 60
$ws.Cells.Item(4, 14).Value = 16000
$ws.Cells.Item(4, 15).Value = 16000
$ws.Cells.Item(4, 16).Value = 16000
$ws.Cells.Item(4, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(4, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 19).Value = 889
$ws.Cells.Item(4, 20).Value = 18
